# Update building block types in the Metabolite extraction template:
#  - "Parameter [bio entity]"        -> "Characteristic [bio entity]"
#  - "Parameter [extraction buffer]" -> "Component [extraction buffer]"
#  - Term Source REF/Accession for bio entity: NFDI4PSO -> EFO term
#  - Remove the now unused "user-specific" value from the extraction buffer
#    term source ref cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2EXT03_Metabolites")

# Header / table column renames (updating the header cell also renames the
# corresponding ListColumn of the annotationTable)
$ws.Range("C1").Value = "Characteristic [bio entity]"
$ws.Range("J1").Value = "Component [extraction buffer]"

# Update the ontology term source ref + accession number for "bio entity"
$ws.Range("D2").Value = "EFO"
$ws.Range("E2").Value = "https://bioregistry.io/EFO:0004964"

# The extraction buffer term source ref value "user-specific" is removed
$ws.Range("K2").Value = ""
